# İş Takip Güncellemesi - 17.12.2025 13:19:11
# Shifts all "İŞ TAKİP" (İş Takip Listesi) İŞE BAŞLAMA / İHALE BİTİŞ dates
# back by one day, plus a couple of DURUM corrections, and shifts the
# matching date columns on the "Güncelleme" sheet back by one day too.

$wb = $excel.ActiveWorkbook

function Shift-CellDateBack1Day {
    param($ws, $row, $col)

    $cell = $ws.Cells.Item($row, $col)
    $oldText = $cell.Value2
    $d = [datetime]::ParseExact($oldText, "yyyy-MM-dd", $null)
    $newText = $d.AddDays(-1).ToString("yyyy-MM-dd")

    # Keep the cell a plain text cell (it was stored as text, not a real
    # date) so the rewritten value round-trips as text too.
    $cell.NumberFormat = "@"
    $cell.Value = $newText
}

# ---------------------------------------------------------------------
# Sheet 1: "İş Takip Listesi"
#   Columns J (İŞE BAŞLAMA/YER TESLİMİ) and K (İHALE BİTİŞ TARİHİ)
#   for rows 2-10 and 33-122 all move back one day.
# ---------------------------------------------------------------------
$wsTakip = $wb.Worksheets.Item("İş Takip Listesi")

$takipRows = @()
$takipRows += 2..10
$takipRows += 33..122

foreach ($r in $takipRows) {
    Shift-CellDateBack1Day $wsTakip $r 10   # J
    Shift-CellDateBack1Day $wsTakip $r 11   # K
}

# DURUM (status) corrections that came along with this update.
$wsTakip.Cells.Item(72, 12).Value = "DEĞERLENDİRMEDE"
$wsTakip.Cells.Item(81, 12).Value = "KESİN ASKIDA"

# ---------------------------------------------------------------------
# Sheet 2: "Güncelleme"
#   Columns I, J, N, P move back one day (F, G, H, K, L, M, O untouched).
# ---------------------------------------------------------------------
$wsGuncelleme = $wb.Worksheets.Item("Güncelleme")

$colIRows = @(5, 7, 9, 11, 19, 24, 29)
$colJRows = @(2, 3, 4, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 25, 27, 28, 29)
$colNRows = @(2, 3, 4, 6, 8, 10, 11, 12, 14, 15, 16, 17, 19, 28, 29)
$colPRows = @(2, 3, 4, 8, 10, 12, 15, 16)

foreach ($r in $colIRows) { Shift-CellDateBack1Day $wsGuncelleme $r 9  }   # I
foreach ($r in $colJRows) { Shift-CellDateBack1Day $wsGuncelleme $r 10 }   # J
foreach ($r in $colNRows) { Shift-CellDateBack1Day $wsGuncelleme $r 14 }   # N
foreach ($r in $colPRows) { Shift-CellDateBack1Day $wsGuncelleme $r 16 }   # P
